$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings formatted with "." as both thousands and
# decimal separators (e.g. "41.655.47"), which Excel would otherwise parse
# as a number. Force the cell to Text format first so the literal string is
# preserved exactly, matching the source data feed.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.655.47"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.474.41"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.72"
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.41"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.10"
$ws.Range("E10").Value = "  +1.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  +8.13%  "

$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.854.87"
$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.61"
$ws.Range("E15").Value = "  -5.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.463.99"
$ws.Range("E16").Value = "  -0.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  +2.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.602.38"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("E19").Value = "  -1.03%  "

$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.38"
$ws.Range("E21").Value = "  -1.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.33"
$ws.Range("E22").Value = "  +1.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.20"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.74"
$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.64"
$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("E28").Value = "  +2.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.86"
$ws.Range("E29").Value = "  +1.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.33"
$ws.Range("E30").Value = "  +1.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.96"
$ws.Range("E31").Value = "  +2.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.51"
$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.59"
$ws.Range("E34").Value = "  +0.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0766"
$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.37"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  +1.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.92"
$ws.Range("E38").Value = "  +0.81%  "

$ws.Range("E39").Value = "  +1.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  -1.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("E41").Value = "  -1.41%  "

$ws.Range("E42").Value = "  +2.72%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.988.87"
$ws.Range("E43").Value = "  +1.36%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0285"
$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.88"
$ws.Range("E45").Value = "  -0.52%  "

$ws.Range("E46").Value = "  +1.22%  "

$ws.Range("E47").Value = "  +3.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.710.98"
$ws.Range("E48").Value = "  -0.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.42"
$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.27"
$ws.Range("E50").Value = "  +2.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.29"
$ws.Range("E51").Value = "  -0.93%  "
